# edit.ps1 - applies the "quantum physics" -> "chemistry" rewrite described by the diff.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: insert a brand new run of text right after a given character
# position, with explicit run formatting (so the run gets its own <w:rPr>
# instead of silently inheriting the paragraph mark's formatting).
# Returns the character position right after the newly inserted text.
# ---------------------------------------------------------------------------
function Insert-FormattedRun {
    param(
        $Doc,
        [int]$Pos,
        [string]$Text,
        [string]$FontName = "Aptos",
        $FontSize = 12,
        $FontColor = 0
    )
    $ins = $Doc.Range($Pos, $Pos)
    $ins.InsertAfter($Text)
    $newPos = $Pos + $Text.Length
    $fmt = $Doc.Range($Pos, $newPos)
    $fmt.Font.Name = $FontName
    $fmt.Font.Size = $FontSize
    $fmt.Font.Color = $FontColor
    return $newPos
}

function Replace-Text {
    param($Doc, [string]$OldText, [string]$NewText)
    $rng = $Doc.Content
    $ok = $rng.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, $NewText, 2)
    if (-not $ok) {
        Write-Host "WARNING: could not find text: $OldText"
    }
    return $ok
}

function Find-EndOf {
    param($Doc, [string]$Text)
    $rng = $Doc.Content
    $ok = $rng.Find.Execute($Text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Host "WARNING: could not find text (for position): $Text"
    }
    return $rng.End
}

# ---------------------------------------------------------------------------
# Title
# ---------------------------------------------------------------------------
Replace-Text $d "Unraveling the Quantum Landscape" "The Ethereal Realm of Chemistry: A Journey into the Intricacies of Matter" | Out-Null

# ---------------------------------------------------------------------------
# Author name
# ---------------------------------------------------------------------------
Replace-Text $d " Ethan Winters" " Eleanor Thompson" | Out-Null

# ---------------------------------------------------------------------------
# Email address: winterse@quantumphysics . edu  ->  eleanor . thompson@scholarlyedu . org
# ---------------------------------------------------------------------------
Replace-Text $d "winterse@quantumphysics" "eleanor" | Out-Null
Replace-Text $d "edu" "thompson@scholarlyedu" | Out-Null

$pos = Find-EndOf $d "thompson@scholarlyedu"
$pos = Insert-FormattedRun $d $pos "." "Aptos" 16 0
$pos = Insert-FormattedRun $d $pos "org" "Aptos" 16 0

# ---------------------------------------------------------------------------
# Body paragraph 1 (three original sentences -> six sentences)
# ---------------------------------------------------------------------------
Replace-Text $d "The realm of quantum physics, a perplexing universe of interconnectedness and uncertainty, continues to captivate scientists and philosophers alike" "In the vast tapestry of science, chemistry stands as a beacon of discovery, illuminating the enigmatic world of matter" | Out-Null

Replace-Text $d " As the foundation of modern physics, quantum mechanics has illuminated the infinitesimal world's intricate workings, revealing a realm where particles dance in probabilistic waves and probabilities govern the fabric of reality" " As a high school teacher, I aim to unveil the secrets of this captivating subject, guiding my students on an enthralling journey through the intricate dance of elements and molecules" | Out-Null

Replace-Text $d " Delving into the quantum landscape signifies a journey into the profound depths of matter, energy, and the fundamental laws that weave the tapestry of existence" " Chemistry is a symphony of reactions and transformations, a saga of particles interacting in a harmonious ballet of creation and destruction" | Out-Null

$pos = Find-EndOf $d " Chemistry is a symphony of reactions and transformations, a saga of particles interacting in a harmonious ballet of creation and destruction"
$pos = Insert-FormattedRun $d $pos "." "Aptos" 12 0
$pos = Insert-FormattedRun $d $pos " Each element, with its unique properties, weaves its magic, forming the very fabric of our universe" "Aptos" 12 0
$pos = Insert-FormattedRun $d $pos "." "Aptos" 12 0
$pos = Insert-FormattedRun $d $pos " It is a realm where the mysteries of the material world unravel, revealing the fundamental forces that shape our lives" "Aptos" 12 0
$pos = Insert-FormattedRun $d $pos "." "Aptos" 12 0
$pos = Insert-FormattedRun $d $pos " It is a subject that combines intellectual rigor with endless fascination, a testament to the boundless curiosity of the human spirit" "Aptos" 12 0

# ---------------------------------------------------------------------------
# Body paragraph 2 (after the first double <w:br/>)
# ---------------------------------------------------------------------------
Replace-Text $d "In this multifaceted domain, particles exhibit both wave-like and particle-like characteristics, a duality that defies classical intuition" "This captivating field unveils the secrets of how matter is composed, how it changes, and how it interacts with its surroundings" | Out-Null

Replace-Text $d " The enigmatic nature of quantum entanglement further astounds, as particles separated by vast distances remain bound in an inseparable embrace of shared fate" " Unraveling the intricacies of these interactions, we glimpse the profound elegance of the natural world" | Out-Null

Replace-Text $d " As we delve deeper into this microscopic realm, the uncertainty principle unveils the inherent interconnectedness of measurements, revealing the inherent limitations of our knowledge" " Chemistry is a gateway to understanding the very essence of things, from the smallest atoms to the grandest molecules" | Out-Null

Replace-Text $d " The quantum world is a tapestry of phenomena that challenge our most fundamental assumptions about reality, inviting us to rethink the very nature of space, time, and existence itself" " As we delve deeper into this realm of discovery, we learn to harness the power of chemical reactions to create new materials, devise novel medicines, and address some of the world's most pressing challenges" | Out-Null

# ---------------------------------------------------------------------------
# Body paragraph 3 (after the second double <w:br/>)
# ---------------------------------------------------------------------------
Replace-Text $d "The study of quantum physics is an endeavor punctuated by both awe and perplexity" "The study of chemistry is not merely an academic pursuit; it is an endeavor that connects us to the world around us" | Out-Null

Replace-Text $d " Its profound implications have rippled across numerous fields, from computation and communication to cosmology and biology" " By understanding the fundamental principles that govern chemical processes, we gain insights into a myriad of phenomena, from the vibrant colors of flowers to the intricate workings of our own bodies" | Out-Null

Replace-Text $d " As scientists continue to unravel the enigmatic tapestry of quantum mechanics, we glimpse the potential for transformative technologies, insights into the nature of consciousness, and a deeper understanding of the universe's fundamental nature" " Chemistry empowers us to make informed decisions about our health, our environment, and our future" | Out-Null

$pos = Find-EndOf $d " Chemistry empowers us to make informed decisions about our health, our environment, and our future"
$pos = Insert-FormattedRun $d $pos "." "Aptos" 12 0
$pos = Insert-FormattedRun $d $pos " It is a discipline that fosters critical thinking, problem-solving skills, and a deep appreciation for the natural world" "Aptos" 12 0
$pos = Insert-FormattedRun $d $pos "." "Aptos" 12 0
$pos = Insert-FormattedRun $d $pos " It prepares us to navigate an increasingly complex world where scientific literacy is essential for informed citizenship" "Aptos" 12 0

# ---------------------------------------------------------------------------
# Summary section
# ---------------------------------------------------------------------------
Replace-Text $d "The realm of quantum physics presents a mind-boggling universe of interconnectedness and uncertainty, challenging classical notions of reality" "Chemistry, an alluring field of scientific exploration, unveils the mysteries of matter, its composition, and its interactions" | Out-Null

Replace-Text $d " Particles exhibit wave-like and particle-like behaviors, entanglements defy distance, and the uncertainty principle reveals the interconnectedness of measurements" " Through the study of chemistry, we gain profound insights into the fundamental forces that shape our universe and the intricate workings of the natural world" | Out-Null

# This sentence spans what used to be two runs split by a <w:lastRenderedPageBreak/>;
# replacing the whole phrase collapses it back into ordinary text.
Replace-Text $d " Quantum physics has profound implications across diverse fields, inspiring new technologies and reshaping our understanding of the cosmos and consciousness" " It empowers us with knowledge and skills essential for navigating an increasingly complex world, enabling us to address global challenges and make informed decisions about our health, our environment, and our future" | Out-Null

Replace-Text $d " As scientists continue to unravel the quantum landscape, the possibilities for transformative advancements are boundless" " Chemistry is a gateway to understanding the very essence of things and fosters a deep appreciation for the harmonious ballet of elements and molecules that comprise our existence" | Out-Null

# ---------------------------------------------------------------------------
# New empty paragraph at the very end of the document body.
# ---------------------------------------------------------------------------
$endRng = $d.Content
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# Let the engine recompute pagination so lastRenderedPageBreak markers
# reflect the new, longer content (the "Summary" heading now starts a new
# page once the document has grown).
# ---------------------------------------------------------------------------
$d.Repaginate()
